$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 1.42
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 8.5
$ws.Range("J2").Value = 1.95
$ws.Range("K2").Value = 2.2
$ws.Range("L2").Value = 8.5
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.67
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("X2").Value = 5.5
$ws.Range("Z2").Value = 8.5
$ws.Range("AA2").Value = 13
$ws.Range("AC2").Value = 8.5
$ws.Range("AD2").Value = 8.5
$ws.Range("AE2").Value = 26
$ws.Range("AF2").Value = 101
$ws.Range("AG2").Value = 15
$ws.Range("AH2").Value = 41
$ws.Range("AI2").Value = 26
$ws.Range("AJ2").Value = 101
$ws.Range("AK2").Value = 67
$ws.Range("AM2").Value = 3.1
$ws.Range("AN2").Value = 7
$ws.Range("AP2").Value = 21
$ws.Range("AT2").Value = 11
$ws.Range("AV2").Value = 9
$ws.Range("AX2").Value = 51
$ws.Range("AY2").Value = 251
$ws.Range("AZ2").Value = 251
$ws.Range("BB2").Value = 126

# Row 3
$ws.Range("G3").Value = 1.9
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.63
$ws.Range("Q3").Value = 2.4
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("X3").Value = 7.5
$ws.Range("Z3").Value = 15
$ws.Range("AE3").Value = 21
$ws.Range("AM3").Value = 3.6
$ws.Range("AW3").Value = 29

# Row 4
$ws.Range("H4").Value = 2.75
$ws.Range("I4").Value = 3
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("O4").Value = 1.67
$ws.Range("P4").Value = 2.1
$ws.Range("Q4").Value = 3.1
$ws.Range("R4").Value = 1.36
$ws.Range("S4").Value = 1.73
$ws.Range("T4").Value = 2
$ws.Range("AC4").Value = 5
$ws.Range("AE4").Value = 21
$ws.Range("AT4").Value = 10
$ws.Range("AY4").Value = 67
$ws.Range("AZ4").Value = 126
